# fix non integer in care home data
# Update the "relpath" column (D) for the "no care-homes" (noCH) rows so
# that the files are looked up under a shared "carehomes" folder instead
# of each study's own folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D19").Value = "data/derived/carehomes/DNK1_agebands_noCH.RDS"
$ws.Range("D20").Value = "data/derived/carehomes/ESP1-2_agebands_noCH.RDS"
$ws.Range("D21").Value = "data/derived/carehomes/GBR3_agebands_noCH.RDS"
$ws.Range("D22").Value = "data/derived/carehomes/CHE1_agebands_noCH.RDS"
$ws.Range("D23").Value = "data/derived/carehomes/CHE2_agebands_noCH.RDS"
$ws.Range("D24").Value = "data/derived/carehomes/NYC_NY_1_agebands_noCH.RDS"

# Move the active selection to match the author's final cursor position.
$ws.Range("D25").Select()
